# Scraper update: insert "height" and "weight" columns immediately
# before the existing "fantasy points" column.
#
# Original layout:  A=index  B=rec_yds  C=rec_td  D=fumbles  E=fantasy points
# New layout:        A=index  B=rec_yds  C=rec_td  D=fumbles  E=height  F=weight  G=fantasy points

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank columns at E:F, pushing the existing "fantasy points"
# column (and its data) from E to G. The inserted cells inherit the
# formatting of the columns around them, same as Excel's normal
# column-insert behaviour.
$ws.Range("E1:F1").EntireColumn.Insert()

# New "height" column/header in E.
$ws.Range("E1").Value = "height"
for ($r = 2; $r -le 6; $r++) {
    $ws.Cells.Item($r, 5).Value = 6.416666666666667
}

# New "weight" column/header in F.
$ws.Range("F1").Value = "weight"
for ($r = 2; $r -le 6; $r++) {
    $ws.Cells.Item($r, 6).Value = 245
}
